# Update the "Förändrad" (changed) date column (C) for rows 2-18
# from 45185 (2023-09-16) to 45204 (2023-10-05), matching the author's
# automatic update of the tracked logging-notification dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value = 45204
    }
}
